$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 23 with the git rm --cached tip
$ws.Range("A23").Value = "git rm .env --cached"
$ws.Range("B23").Value = "esto se hace cuando en el .gitignore quiero ignorar un archivo que ya subi anteriormente, reemplazo el .env con el nombre del archivo, luego hago el .add commit y push"

# Match formatting used by the other wrapped-text rows (e.g. row 17)
$ws.Range("A23").VerticalAlignment = -4108
$ws.Range("B23").WrapText = $true
$ws.Rows.Item(23).RowHeight = 45

# Update selection / view state to the newly active cell
$ws.Activate()
$ws.Range("B23").Select()

# Page setup was touched (portrait orientation) in the edited workbook
$ws.PageSetup.Orientation = 1
